$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.916.93"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "1.642.75"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.78"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5056"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2581"
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06405"
$ws.Range("E9").Value = "  -0.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.62"
$ws.Range("E10").Value = "  +0.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07778"
$ws.Range("E11").Value = "  +0.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.290"
$ws.Range("E12").Value = "  +0.97%  "
$ws.Range("D13").Value = "1.653.87"
$ws.Range("E13").Value = "  +0.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5441"
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("D15").Value = "0.0₅7886"
$ws.Range("E15").Value = "  -0.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.99"
$ws.Range("E16").Value = "  +2.49%  "
$ws.Range("D17").Value = "25.966.85"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.005"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "197.88"
$ws.Range("E19").Value = "  -2.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.417"
$ws.Range("E20").Value = "  +2.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.969"
$ws.Range("E21").Value = "  -0.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.008"
$ws.Range("E22").Value = "  +0.50%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.868"
$ws.Range("E24").Value = "  -4.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "140.88"
$ws.Range("E25").Value = "  -0.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1144"
$ws.Range("E26").Value = "  -0.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.886"
$ws.Range("E27").Value = "  +2.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.75"
$ws.Range("E28").Value = "  +0.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.246"
$ws.Range("E29").Value = "  +0.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05027"
$ws.Range("E30").Value = "  -0.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.273"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.210"
$ws.Range("E32").Value = "  +0.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.535"
$ws.Range("E33").Value = "  -0.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.378"
$ws.Range("E34").Value = "  +1.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.8948"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.614"
$ws.Range("E36").Value = "  -0.95%  "
$ws.Range("D37").Value = "1.146.02"
$ws.Range("E37").Value = "  -0.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5561"
$ws.Range("E38").Value = "  -1.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01565"
$ws.Range("E39").Value = "  -0.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.008"
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.698"
$ws.Range("E41").Value = "  +0.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8204"
$ws.Range("E42").Value = "  +1.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.08"
$ws.Range("E43").Value = "  +0.48%  "
$ws.Range("E44").Value = "  +8.32%  "
$ws.Range("D45").Value = "1.782.73"
$ws.Range("E45").Value = "  +0.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4529"
$ws.Range("E46").Value = "  +0.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.41"
$ws.Range("E47").Value = "  +0.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.008"
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05076"
$ws.Range("E49").Value = "  +0.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.007"
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.09556"
$ws.Range("E51").Value = "  +2.79%  "
